$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    $origStyle = $range.Style
    $range.Value = "'" + $text
    $range.Style = $origStyle
}

$v0 = '58.722.03'
Set-TextValue $ws.Range("D2") $v0

$v1 = '  -2.25%  '
$ws.Range("E2").Value = $v1

$v2 = '2.308.87'
Set-TextValue $ws.Range("D3") $v2

$v3 = '  -4.72%  '
$ws.Range("E3").Value = $v3

$v4 = '  -0.09%  '
$ws.Range("E4").Value = $v4

$v5 = '551.06'
Set-TextValue $ws.Range("D5") $v5

$v6 = '  -0.56%  '
$ws.Range("E5").Value = $v6

$v7 = '131.23'
Set-TextValue $ws.Range("D6") $v7

$v8 = '  -4.60%  '
$ws.Range("E6").Value = $v8

$v9 = '  -0.03%  '
$ws.Range("E7").Value = $v9

$v10 = '0.569'
Set-TextValue $ws.Range("D8") $v10

$v11 = '  -4.12%  '
$ws.Range("E8").Value = $v11

$v12 = '0.102'
Set-TextValue $ws.Range("D9") $v12

$v13 = '  -3.14%  '
$ws.Range("E9").Value = $v13

$v14 = '5.55'
Set-TextValue $ws.Range("D10") $v14

$v15 = '  -2.15%  '
$ws.Range("E10").Value = $v15

$v16 = '  +0.73%  '
$ws.Range("E11").Value = $v16

$v17 = '0.337'
Set-TextValue $ws.Range("D12") $v17

$v18 = '  -5.04%  '
$ws.Range("E12").Value = $v18

$v19 = '23.80'
Set-TextValue $ws.Range("D13") $v19

$v20 = '  -5.95%  '
$ws.Range("E13").Value = $v20

$v21 = '2.720.16'
Set-TextValue $ws.Range("D14") $v21

$v22 = '  -4.77%  '
$ws.Range("E14").Value = $v22

$v23 = '58.751.60'
Set-TextValue $ws.Range("D15") $v23

$v24 = '  -2.10%  '
$ws.Range("E15").Value = $v24

$v25 = '0.0000133'
Set-TextValue $ws.Range("D16") $v25

$v26 = '  -3.18%  '
$ws.Range("E16").Value = $v26

$v27 = '2.306.67'
Set-TextValue $ws.Range("D17") $v27

$v28 = '  -4.89%  '
$ws.Range("E17").Value = $v28

$v29 = '10.77'
Set-TextValue $ws.Range("D18") $v29

$v30 = '  -4.76%  '
$ws.Range("E18").Value = $v30

$v31 = '4.37'
Set-TextValue $ws.Range("D19") $v31

$v32 = '  -1.98%  '
$ws.Range("E19").Value = $v32

$v33 = '315.56'
Set-TextValue $ws.Range("D20") $v33

$v34 = '  -4.02%  '
$ws.Range("E20").Value = $v34

$v35 = '6.53'
Set-TextValue $ws.Range("D21") $v35

$v36 = '  -2.39%  '
$ws.Range("E21").Value = $v36

$v37 = '  -0.06%  '
$ws.Range("E22").Value = $v37

$v38 = '63.07'
Set-TextValue $ws.Range("D23") $v38

$v39 = '  -4.49%  '
$ws.Range("E23").Value = $v39

$v40 = '0.171'
Set-TextValue $ws.Range("D24") $v40

$v41 = '  -4.34%  '
$ws.Range("E24").Value = $v41

$v42 = '  +0.08%  '
$ws.Range("E25").Value = $v42

$v43 = '8.23'
Set-TextValue $ws.Range("D26") $v43

$v44 = '  -4.36%  '
$ws.Range("E26").Value = $v44

$v45 = '1.32'
Set-TextValue $ws.Range("D27") $v45

$v46 = '  -6.83%  '
$ws.Range("E27").Value = $v46

$v47 = '1.77'
Set-TextValue $ws.Range("D28") $v47

$v48 = '  -0.08%  '
$ws.Range("E28").Value = $v48

$v49 = '170.95'
Set-TextValue $ws.Range("D29") $v49

$v50 = '  +0.86%  '
$ws.Range("E29").Value = $v50

$v51 = '0.0'
$v51 = $v51 + [string][char]8323
$v51 = $v51 + '0733'
Set-TextValue $ws.Range("D30") $v51

$v52 = '  -5.94%  '
$ws.Range("E30").Value = $v52

$v53 = '5.86'
Set-TextValue $ws.Range("D31") $v53

$v54 = '  -3.77%  '
$ws.Range("E31").Value = $v54

$v55 = '  +0.01%  '
$ws.Range("E32").Value = $v55

$v56 = '0.388'
Set-TextValue $ws.Range("D33") $v56

$v57 = '  -3.96%  '
$ws.Range("E33").Value = $v57

$v58 = '  +0.04%  '
$ws.Range("E34").Value = $v58

$v59 = '17.78'
Set-TextValue $ws.Range("D35") $v59

$v60 = '  -4.62%  '
$ws.Range("E35").Value = $v60

$v61 = '  +0.01%  '
$ws.Range("E36").Value = $v61

$v62 = '  -3.43%  '
$ws.Range("E37").Value = $v62

$v63 = '4.00'
Set-TextValue $ws.Range("D38") $v63

$v64 = '  -5.09%  '
$ws.Range("E38").Value = $v64

$v65 = 'OKB'
$ws.Range("B39").Value = $v65

$v66 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C39").Value = $v66

$v67 = '38.52'
Set-TextValue $ws.Range("D39") $v67

$v68 = '  -2.26%  '
$ws.Range("E39").Value = $v68

$v69 = 'Stacks'
$ws.Range("B40").Value = $v69

$v70 = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("C40").Value = $v70

$v71 = '1.52'
Set-TextValue $ws.Range("D40") $v71

$v72 = '  -5.28%  '
$ws.Range("E40").Value = $v72

$v73 = '296.66'
Set-TextValue $ws.Range("D41") $v73

$v74 = '  -9.32%  '
$ws.Range("E41").Value = $v74

$v75 = '142.15'
Set-TextValue $ws.Range("D42") $v75

$v76 = '  +0.60%  '
$ws.Range("E42").Value = $v76

$v77 = '3.46'
Set-TextValue $ws.Range("D43") $v77

$v78 = '  -5.55%  '
$ws.Range("E43").Value = $v78

$v79 = '0.0947'
Set-TextValue $ws.Range("D44") $v79

$v80 = '  -2.59%  '
$ws.Range("E44").Value = $v80

$v81 = '0.0499'
Set-TextValue $ws.Range("D45") $v81

$v82 = '  -3.52%  '
$ws.Range("E45").Value = $v82

$v83 = '18.62'
Set-TextValue $ws.Range("D46") $v83

$v84 = '  -5.18%  '
$ws.Range("E46").Value = $v84

$v85 = '0.556'
Set-TextValue $ws.Range("D47") $v85

$v86 = '  -3.43%  '
$ws.Range("E47").Value = $v86

$v87 = '0.0214'
Set-TextValue $ws.Range("D48") $v87

$v88 = '  -4.29%  '
$ws.Range("E48").Value = $v88

$v89 = '11.03'
Set-TextValue $ws.Range("D49") $v89

$v90 = '  -0.21%  '
$ws.Range("E49").Value = $v90

$v91 = '  +0.24%  '
$ws.Range("E50").Value = $v91

$v92 = '  -1.04%  '
$ws.Range("E51").Value = $v92

